$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (shifts existing rows 3-14 down to 4-15)
$ws.Rows.Item(3).Insert()

# Copy the static columns from row 2 (the template record) into the new row 3
$ws.Range("A3:C3").Value = $ws.Range("A2:C2").Value2
$ws.Range("E3:L3").Value = $ws.Range("E2:L2").Value2
$ws.Range("Q3:R3").Value = $ws.Range("Q2:R2").Value2
$ws.Range("T3").Value = $ws.Range("T2").Value2

# Match the date format used by the other Fecha cells
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat

# New record's values
$ws.Range("D3").Value = 44881
$ws.Range("M3").Value = 440
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6500
$ws.Range("S3").Value = 3250
